{"js": "// Replace the date line and the 25 multiplication problems in the table\n// with the values from the target revision. Each \"before\" string is\n// unique within the document, so a direct search-and-replace per pair\n// is unambiguous.\nconst replacements = [\n  [\"2024-04-09 Tuesday\", \"2024-04-10 Wednesday\"],\n  [\"531\u00d79=4779\", \"422\u00d74=1688\"],\n  [\"864\u00d76=5184\", \"478\u00d78=3824\"],\n  [\"415\u00d76=2490\", \"937\u00d76=5622\"],\n  [\"938\u00d72=1876\", \"178\u00d73=534\"],\n  [\"590\u00d75=2950\", \"909\u00d75=4545\"],\n  [\"996\u00d77=6972\", \"105\u00d74=420\"],\n  [\"860\u00d78=6880\", \"111\u00d73=333\"],\n  [\"161\u00d78=1288\", \"143\u00d75=715\"],\n  [\"599\u00d77=4193\", \"654\u00d74=2616\"],\n  [\"692\u00d74=2768\", \"979\u00d78=7832\"],\n  [\"164\u00d73=492\", \"189\u00d76=1134\"],\n  [\"237\u00d74=948\", \"784\u00d75=3920\"],\n  [\"838\u00d74=3352\", \"213\u00d74=852\"],\n  [\"536\u00d74=2144\", \"318\u00d75=1590\"],\n  [\"923\u00d77=6461\", \"429\u00d78=3432\"],\n  [\"837\u00d72=1674\", \"583\u00d74=2332\"],\n  [\"914\u00d72=1828\", \"973\u00d72=1946\"],\n  [\"985\u00d77=6895\", \"186\u00d75=930\"],\n  [\"444\u00d72=888\", \"390\u00d73=1170\"],\n  [\"264\u00d75=1320\", \"595\u00d75=2975\"],\n  [\"501\u00d72=1002\", \"800\u00d73=2400\"],\n  [\"861\u00d72=1722\", \"320\u00d78=2560\"],\n  [\"924\u00d75=4620\", \"348\u00d76=2088\"],\n  [\"624\u00d79=5616\", \"371\u00d76=2226\"],\n  [\"509\u00d73=1527\", \"114\u00d72=228\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${before}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 multiplication problems in the table\n# with the values from the target revision. Each \"before\" string is\n# unique within the document, so a direct Find/Replace per pair is\n# unambiguous. NOTE: the before/after strings are written as literal\n# UTF-8 text (not built with \"+\" concatenation) to avoid the numeric\n# string coercion quirk of this shell on digit-only operands.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-04-09 Tuesday\", \"2024-04-10 Wednesday\"),\n    @(\"531\u00d79=4779\", \"422\u00d74=1688\"),\n    @(\"864\u00d76=5184\", \"478\u00d78=3824\"),\n    @(\"415\u00d76=2490\", \"937\u00d76=5622\"),\n    @(\"938\u00d72=1876\", \"178\u00d73=534\"),\n    @(\"590\u00d75=2950\", \"909\u00d75=4545\"),\n    @(\"996\u00d77=6972\", \"105\u00d74=420\"),\n    @(\"860\u00d78=6880\", \"111\u00d73=333\"),\n    @(\"161\u00d78=1288\", \"143\u00d75=715\"),\n    @(\"599\u00d77=4193\", \"654\u00d74=2616\"),\n    @(\"692\u00d74=2768\", \"979\u00d78=7832\"),\n    @(\"164\u00d73=492\", \"189\u00d76=1134\"),\n    @(\"237\u00d74=948\", \"784\u00d75=3920\"),\n    @(\"838\u00d74=3352\", \"213\u00d74=852\"),\n    @(\"536\u00d74=2144\", \"318\u00d75=1590\"),\n    @(\"923\u00d77=6461\", \"429\u00d78=3432\"),\n    @(\"837\u00d72=1674\", \"583\u00d74=2332\"),\n    @(\"914\u00d72=1828\", \"973\u00d72=1946\"),\n    @(\"985\u00d77=6895\", \"186\u00d75=930\"),\n    @(\"444\u00d72=888\", \"390\u00d73=1170\"),\n    @(\"264\u00d75=1320\", \"595\u00d75=2975\"),\n    @(\"501\u00d72=1002\", \"800\u00d73=2400\"),\n    @(\"861\u00d72=1722\", \"320\u00d78=2560\"),\n    @(\"924\u00d75=4620\", \"348\u00d76=2088\"),\n    @(\"624\u00d79=5616\", \"371\u00d76=2226\"),\n    @(\"509\u00d73=1527\", \"114\u00d72=228\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
